$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the entire content of column A (removes the title "State-By-State
# Breakdown" from A1 and the "State" header / all state abbreviations from
# A2:A55). Columns B:J are left untouched.
$ws.Columns.Item(1).ClearContents()

# Update the active selection to match the post-edit workbook (A1:C1).
$ws.Range("A1:C1").Select()
